$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.109663917146452
$ws.Range("D2").Value = 1.111241962773572
$ws.Range("E2").Value = 1.108271293359297
$ws.Range("F2").Value = 1.117314121119132
$ws.Range("I2").Value = 1.031262405943017
$ws.Range("J2").Value = 1.114402818680321
$ws.Range("K2").Value = 1.113837281912491
$ws.Range("L2").Value = 1.110873942118266
$ws.Range("M2").Value = 1.11989458862937
$ws.Range("N2").Value = 1.115985398506637
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.11266461397181
$ws.Range("D3").Value = 1.11399711166576
$ws.Range("E3").Value = 1.110941694722924
$ws.Range("F3").Value = 1.120042146674648
$ws.Range("I3").Value = 1.031350524314284
$ws.Range("J3").Value = 1.117066580033014
$ws.Range("K3").Value = 1.116411878131243
$ws.Range("L3").Value = 1.113363434384528
$ws.Range("M3").Value = 1.122443236909063
$ws.Range("N3").Value = 1.118652942705989
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.114594975957789
$ws.Range("D4").Value = 1.115769031894531
$ws.Range("E4").Value = 1.112658950624541
$ws.Range("F4").Value = 1.121796086414583
$ws.Range("I4").Value = 1.031404225244571
$ws.Range("J4").Value = 1.118779015596356
$ws.Range("K4").Value = 1.11806670561324
$ws.Range("L4").Value = 1.114963352472252
$ws.Range("M4").Value = 1.124080834974883
$ws.Range("N4").Value = 1.12036781012425
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.115403873929366
$ws.Range("D5").Value = 1.116511420320558
$ws.Range("E5").Value = 1.113378397488971
$ws.Range("F5").Value = 1.12253081200276
$ws.Range("I5").Value = 1.031426010716302
$ws.Range("J5").Value = 1.119496311918708
$ws.Range("K5").Value = 1.118759802342718
$ws.Range("L5").Value = 1.115633401261866
$ws.Range("M5").Value = 1.12476658240167
$ws.Range("N5").Value = 1.121086125089654
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.115539539656499
$ws.Range("D6").Value = 1.116635924429014
$ws.Range("E6").Value = 1.113499051853929
$ws.Range("F6").Value = 1.122654023371029
$ws.Range("I6").Value = 1.031429622352942
$ws.Range("J6").Value = 1.119616597960679
$ws.Range("K6").Value = 1.118876026209575
$ws.Range("L6").Value = 1.115745757330273
$ws.Range("M6").Value = 1.124881566043929
$ws.Range("N6").Value = 1.121206581951601
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.114605794709837
$ws.Range("D7").Value = 1.115778961561522
$ws.Range("E7").Value = 1.112668573593524
$ws.Range("F7").Value = 1.121805914099582
$ws.Range("I7").Value = 1.031404519443683
$ws.Range("J7").Value = 1.118788610315459
$ws.Range("K7").Value = 1.118075976905003
$ws.Range("L7").Value = 1.114972315657613
$ws.Range("M7").Value = 1.124090008482761
$ws.Range("N7").Value = 1.120377418468955
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.110680400782044
$ws.Range("D8").Value = 1.112175366131243
$ws.Range("E8").Value = 1.109176018801933
$ws.Range("F8").Value = 1.118238446465147
$ws.Range("I8").Value = 1.031292875218688
$ws.Range("J8").Value = 1.115305409366931
$ws.Range("K8").Value = 1.114709719350563
$ws.Range("L8").Value = 1.111717584238815
$ws.Range("M8").Value = 1.12075834728555
$ws.Range("N8").Value = 1.116889270975551
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.103673395772726
$ws.Range("D9").Value = 1.105739131032289
$ws.Range("E9").Value = 1.102936871926895
$ws.Range("F9").Value = 1.111862622786674
$ws.Range("I9").Value = 1.031070558281898
$ws.Range("J9").Value = 1.109078672524021
$ws.Range("K9").Value = 1.10868982193841
$ws.Range("L9").Value = 1.105895515226617
$ws.Range("M9").Value = 1.11479608745897
$ws.Range("N9").Value = 1.110653691452107
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.098936605961478
$ws.Range("D10").Value = 1.101385807155698
$ws.Range("E10").Value = 1.098716030942434
$ws.Range("F10").Value = 1.107547446775642
$ws.Range("I10").Value = 1.03090488898155
$ws.Range("J10").Value = 1.104863279816671
$ws.Range("K10").Value = 1.104612979713878
$ws.Range("L10").Value = 1.101951565035478
$ws.Range("M10").Value = 1.110755489574353
$ws.Range("N10").Value = 1.106432312403601
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.096868868471015
$ws.Range("D11").Value = 1.099484917606227
$ws.Range("E11").Value = 1.096872797509979
$ws.Range("F11").Value = 1.105662584515994
$ws.Range("I11").Value = 1.030828952181815
$ws.Range("J11").Value = 1.103021715628051
$ws.Range("K11").Value = 1.102831597089646
$ws.Range("L11").Value = 1.100227994154452
$ws.Range("M11").Value = 1.108989296252733
$ws.Range("N11").Value = 1.104588132982602
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.096098214354031
$ws.Range("D12").Value = 1.098776368096936
$ws.Range("E12").Value = 1.096185710337711
$ws.Range("F12").Value = 1.104959915621905
$ws.Range("I12").Value = 1.030800109420914
$ws.Range("J12").Value = 1.102335143092599
$ws.Range("K12").Value = 1.102167409723844
$ws.Range("L12").Value = 1.099585323562275
$ws.Range("M12").Value = 1.108330676094598
$ws.Range("N12").Value = 1.103900585435399
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.096263641748833
$ws.Range("D13").Value = 1.098928467819781
$ws.Range("E13").Value = 1.096333204177667
$ws.Range("F13").Value = 1.105110757201204
$ws.Range("I13").Value = 1.030806325179219
$ws.Range("J13").Value = 1.10248253133915
$ws.Range("K13").Value = 1.102309994824835
$ws.Range("L13").Value = 1.099723291258241
$ws.Range("M13").Value = 1.108472070377202
$ws.Range("N13").Value = 1.104048182990166
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.096805219671215
$ws.Range("D14").Value = 1.099426399690961
$ws.Range("E14").Value = 1.096816052718061
$ws.Range("F14").Value = 1.105604554169286
$ws.Range("I14").Value = 1.030826581045043
$ws.Range("J14").Value = 1.102965015536486
$ws.Range("K14").Value = 1.102776746713346
$ws.Range("L14").Value = 1.100174921470785
$ws.Range("M14").Value = 1.108934907582771
$ws.Range("N14").Value = 1.104531352370403
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.097138555598367
$ws.Range("D15").Value = 1.099732861324502
$ws.Range("E15").Value = 1.097113227229177
$ws.Range("F15").Value = 1.105908458562147
$ws.Range("I15").Value = 1.030838976854376
$ws.Range("J15").Value = 1.103261951544903
$ws.Range("K15").Value = 1.103063993726481
$ws.Range("L15").Value = 1.100452857285128
$ws.Range("M15").Value = 1.109219732753995
$ws.Range("N15").Value = 1.104828710062012
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.09907347515443
$ws.Range("D16").Value = 1.101511620913149
$ws.Range("E16").Value = 1.098838024594939
$ws.Range("F16").Value = 1.107672186606111
$ws.Range("I16").Value = 1.030909839699117
$ws.Range("J16").Value = 1.104985148003806
$ws.Range("K16").Value = 1.104730858050522
$ws.Range("L16").Value = 1.102065612473671
$ws.Range("M16").Value = 1.110872349302365
$ws.Range("N16").Value = 1.106554353657539
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.100282665340999
$ws.Range("D17").Value = 1.102623077760878
$ws.Range("E17").Value = 1.09991571192859
$ws.Range("F17").Value = 1.108774083902495
$ws.Range("I17").Value = 1.030953161692707
$ws.Range("J17").Value = 1.106061645209735
$ws.Range("K17").Value = 1.105772071939199
$ws.Range("L17").Value = 1.103072958834656
$ws.Range("M17").Value = 1.111904491543604
$ws.Range("N17").Value = 1.107632379612977
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.100986363312422
$ws.Range("D18").Value = 1.103269846550508
$ws.Range("E18").Value = 1.10054281134174
$ws.Range("F18").Value = 1.109415229316852
$ws.Range("I18").Value = 1.03097802572104
$ws.Range("J18").Value = 1.106687983996954
$ws.Range("K18").Value = 1.106377848541281
$ws.Range("L18").Value = 1.103659006822066
$ws.Range("M18").Value = 1.11250492803532
$ws.Range("N18").Value = 1.108259607873116
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.101226037333383
$ws.Range("D19").Value = 1.103490121987659
$ws.Range("E19").Value = 1.10075638466616
$ws.Range("F19").Value = 1.109633579245938
$ws.Range("I19").Value = 1.030986435182647
$ws.Range("J19").Value = 1.106901286749591
$ws.Range("K19").Value = 1.106584143086504
$ws.Range("L19").Value = 1.103858578709032
$ws.Range("M19").Value = 1.112709393436663
$ws.Range("N19").Value = 1.108473213540127
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.100153097108677
$ws.Range("D20").Value = 1.102503987360159
$ws.Range("E20").Value = 1.099800241781375
$ws.Range("F20").Value = 1.108656024051127
$ws.Range("I20").Value = 1.030948555577378
$ws.Range("J20").Value = 1.105946309637092
$ws.Range("K20").Value = 1.105660520053356
$ws.Range("L20").Value = 1.102965037933165
$ws.Range("M20").Value = 1.111793917996838
$ws.Range("N20").Value = 1.107516880250592
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.096645811171409
$ws.Range("D21").Value = 1.099279840207658
$ws.Range("E21").Value = 1.096673933599958
$ws.Range("F21").Value = 1.105459214307723
$ws.Range("I21").Value = 1.030820633810166
$ws.Range("J21").Value = 1.102823006632003
$ws.Range("K21").Value = 1.102639369701784
$ws.Range("L21").Value = 1.100041996210248
$ws.Range("M21").Value = 1.108798685403352
$ws.Range("N21").Value = 1.104389141796987
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.09442552246328
$ws.Range("D22").Value = 1.097238327995991
$ws.Range("E22").Value = 1.094694205827996
$ws.Range("F22").Value = 1.103434470552045
$ws.Range("I22").Value = 1.030736519240048
$ws.Range("J22").Value = 1.100844557651044
$ws.Range("K22").Value = 1.100725329258794
$ws.Range("L22").Value = 1.098189890184713
$ws.Range("M22").Value = 1.10690050839386
$ws.Range("N22").Value = 1.102407883191544
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.095604006321375
$ws.Range("D23").Value = 1.098321964344612
$ws.Range("E23").Value = 1.095745062640598
$ws.Range("F23").Value = 1.104509257019102
$ws.Range("I23").Value = 1.030781461088725
$ws.Range("J23").Value = 1.101894795034212
$ws.Range("K23").Value = 1.101741404346822
$ws.Range("L23").Value = 1.099173108131329
$ws.Range("M23").Value = 1.107908214574475
$ws.Range("N23").Value = 1.10345961203226
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.100211648361081
$ws.Range("D24").Value = 1.102557803894004
$ws.Range("E24").Value = 1.099852422392712
$ws.Range("F24").Value = 1.108709375070789
$ws.Range("I24").Value = 1.030950638133713
$ws.Range("J24").Value = 1.105998429645169
$ws.Range("K24").Value = 1.105710930320843
$ws.Range("L24").Value = 1.103013807435717
$ws.Range("M24").Value = 1.111843886345795
$ws.Range("N24").Value = 1.107569074275059
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.105496039281522
$ws.Range("D25").Value = 1.107413730197124
$ws.Range("E25").Value = 1.104560336963234
$ws.Range("F25").Value = 1.113521983335021
$ws.Range("I25").Value = 1.03113109035066
$ws.Range("J25").Value = 1.110699420692694
$ws.Range("K25").Value = 1.110256991706428
$ws.Range("L25").Value = 1.10741137386605
$ws.Range("M25").Value = 1.112276741269082
$ws.Range("N25").Value = 1.112276741269082
